## Updated cryptos list on Fri Apr 21 07:43:10 UTC 2023 with GitHub Actions
## Refresh the "Price" (D) and "Volume(1h)" (E) columns for rows 2-51 with the
## latest scraped values. Price cells are forced to text ("@" number format,
## then the style is reset to "Normal" so no stray style index is left behind)
## because several price strings look numeric (e.g. "1.010", "0.00001069",
## "12.80") and a plain .Value assignment would otherwise have Excel coerce
## them into numbers/scientific notation and silently drop the formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.199.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.930.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4744"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4049"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08506"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.946.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("E14").Value = "  -3.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.119"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.011"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06606"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.775"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.235.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.178.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.165"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.770"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9805"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09598"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E34").Value = "  -4.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.664"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.587"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.302"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02322"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06177"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.237"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6188"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1909"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.323"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5905"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.045"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.397"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06771"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("E51").Value = "  -2.23%  "
